$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Parameter name changes from unit_incremental_heat_rate to unit_flow_coefficient
$ws.Range("B5").Value = "unit_flow_coefficient"

# Header row: relationship class changes from unit__node__node to unit__from_node__user_constraint
$ws.Range("B1").Value = "unit__from_node__user_constraint"

# Replace formulas with plain numeric values (inverse of the heat rate -> flow coefficient)
$ws.Range("B6").Value = 0.75
$ws.Range("B7").Value = 0.7
$ws.Range("B8").Value = 0.65

# Update the active selection to match the saved workbook state
$ws.Range("B4").Select()
